# Updated capital structure database
# - Remove the old last data row (row 6, Intercorp Financial Services Inc.)
# - Shuffle company names/order for rows 2-5 and refresh all of their metrics

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop row 6 entirely (worksheet shrinks from A1:AQ6 to A1:AQ5)
$ws.Rows.Item(6).Delete()

# 2) Company names (column B) for rows 2-5
# B2 keeps a purely numeric-looking label ("3"), so force Text format first
# so it round-trips as a string instead of becoming the number 3.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"

$ws.Range("B3").Value = "Banco de Crédito del Perú S.A. (BVL:CREDITC1)"
$ws.Range("B4").Value = "Banco BBVA Perú, S.A. (BVL:BBVAC1)"
$ws.Range("B5").Value = "Intercorp Financial Services Inc. (BVL:IFS)"

# 3) Refreshed numeric metrics for rows 2-5
$values = @{
  "D2"  = -0.0133
  "E2"  = -0.0906
  "K2"  = 820.6
  "L2"  = 0.2220478406753978
  "M2"  = 861.255
  "N2"  = 0.04503600244722517
  "O2"  = 1.049543017304411
  "P2"  = 861.2
  "Q2"  = 0.04503312643473805
  "R2"  = 1.049475993175725
  "S2"  = 0.05500000000000682
  "T2"  = 0.00006386029689233365
  "U2"  = 18967.9
  "V2"  = 0.9918530409910217
  "W2"  = 0.08342842178117811
  "X2"  = 0.06681065503914499
  "Y2"  = 0.01661776674203312
  "Z2"  = 0.2572061713891193
  "AB2" = 0.04196599435506897
  "AC2" = -0.04196599435506897
  "AD2" = 25050.4
  "AF2" = 25050.4
  "AG2" = 6082.5
  "AH2" = 0.5670834267138436
  "AI2" = 0.7150004994933711
  "AJ2" = 0.2413096777776896
  "AK2" = 0.3785568473200727

  "D3"  = -0.0103
  "E3"  = -0.0906
  "K3"  = 446
  "L3"  = 0.2134277647509212
  "M3"  = 510.6
  "N3"  = 0.04389312977099237
  "O3"  = 1.144843049327354
  "P3"  = 510.6
  "Q3"  = 0.04389312977099237
  "R3"  = 1.144843049327354
  "U3"  = 8483.799999999999
  "V3"  = 0.7292999105976206
  "W3"  = 0.08342842178117811
  "X3"  = 0.05939141747663106
  "Y3"  = 0.02403700430454705
  "Z3"  = 0.2742496423743717
  "AB3" = 0.04115707489571348
  "AC3" = -0.04115707489571348
  "AD3" = 12816.6
  "AF3" = 12816.6
  "AG3" = 4332.800000000001
  "AH3" = 0.5242091830474367
  "AI3" = 0.7153038615447295
  "AJ3" = 0.2713834744701108
  "AK3" = 0.4592798312468863

  "D4"  = -0.0133
  "E4"  = -0.078
  "K4"  = 258.5
  "L4"  = 0.2807646356033452
  "M4"  = 156.4
  "N4"  = 0.0399081398315897
  "O4"  = 0.6050290135396519
  "P4"  = 156.4
  "Q4"  = 0.0399081398315897
  "R4"  = 0.6050290135396519
  "U4"  = 5776
  "V4"  = 1.473845368716509
  "W4"  = 0.09930849020361122
  "X4"  = 0.07281762749817698
  "Y4"  = 0.02649086270543424
  "Z4"  = 0.2768239956222902
  "AB4" = 0.04196599435506897
  "AC4" = -0.04196599435506897
  "AD4" = 6965.1
  "AF4" = 6965.1
  "AG4" = 1189.1
  "AH4" = 0.6399334809492746
  "AI4" = 0.7307607566648832
  "AJ4" = 0.2327871419901725
  "AK4" = 0.3166458072590739

  "D5"  = -0.0549
  "E5"  = -0.2
  "K5"  = 116.1
  "L5"  = 0.1694395796847636
  "M5"  = 194.255
  "N5"  = 0.05438422128279067
  "O5"  = 1.673169681309216
  "P5"  = 194.2
  "Q5"  = 0.05436882331532238
  "R5"  = 1.672695951765719
  "S5"  = 0.05500000000000682
  "T5"  = 0.0002831329952897317
  "U5"  = 4708.1
  "V5"  = 1.31809401159047
  "W5"  = 0.04661340185490023
  "X5"  = 0.06681065503914499
  "Y5"  = -0.02019725318424476
  "Z5"  = 0.2001986793665634
  "AB5" = 0.04422161473855393
  "AC5" = -0.04422161473855393
  "AD5" = 5268.7
  "AF5" = 5268.7
  "AG5" = 560.5999999999995
  "AH5" = 0.5959663371264393
  "AI5" = 0.6944836222236868
  "AJ5" = 0.1356563823351481
  "AK5" = 0.1947609783212894
}

foreach ($ref in $values.Keys) {
  $ws.Range($ref).Value = $values[$ref]
}
